$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 46): Amon / o3 missing-variable entry for the
# pextra identified-missing-variables list.

# Column C ("prio") holds text like "1" elsewhere in the sheet (shared
# string, not a number) - force text formatting before entry so it is
# stored as a string rather than being auto-coerced to a numeric value.
$ws.Range("C46").NumberFormat = "@"

$ws.Range("A46").Value = "Amon"
$ws.Range("B46").Value = "o3"
$ws.Range("C46").Value = "1"
$ws.Range("D46").Value = "longitude latitude plev19 time"
$ws.Range("E46").Value = "Mole Fraction of O3"
$ws.Range("F46").Value = "mol mol-1"
$ws.Range("G46").Formula = '=HYPERLINK("http://clipc-services.ceda.ac.uk/dreq/u/1d4594c97188efd47935238a429e02e4.html","web")'
$ws.Range("H46").Value = "tm5 code name = o3|ifs code name = 203.128"
$ws.Range("I46").Value = "automatic"
$ws.Range("J46").Value = "Mole fraction is used in the construction mole_fraction_of_X_in_Y, where X is a material constituent of Y."
$ws.Range("K46").Value = "AerChemMIP,C4MIP,CFMIP,CMIP,DAMIP,FAFMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,LUMIP,RFMIP,VolMIP"

$ws.Range("C46").NumberFormat = "General"

# Match the row height used by the new row.
$ws.Rows.Item(46).RowHeight = 15

# Move the view so the newly added row is visible/selected, as in the
# saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Rows.Item(46).Select()
